# cg: bsa - update district list.
#
# The "choices" sheet's region_list / district_list choice lists are
# replaced with an updated (and corrected) list of Congo regions and
# districts. Previously there were 3 regions (Bouenza, Niari, Pool) with
# 17 districts; now there are 7 regions with 37 districts, and a couple
# of district names get accent corrections (Goma Tse-Tse -> Goma
# Tsé-Tsé, Igne-Ngabe-Mayama -> Igné-Ngabe-Mayama, Loutete -> Loutété).
#
# The "choices" sheet also becomes the active/selected tab (it was
# "survey" before).

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- New region_list (rows 15-21) ---------------------------------------
$regions = @(
    'Bouenza',
    'Brazzaville',
    'Cuvette',
    'Kouilou',
    'Niari',
    'Plateaux',
    'Pool'
)

# --- New district_list (rows 22-58): name, region ------------------------
$districts = @(
    @('Loudima', 'Bouenza'),
    @('Loutété', 'Bouenza'),
    @('Madingou', 'Bouenza'),
    @('Mouyondzi', 'Bouenza'),
    @('Nkayi', 'Bouenza'),
    @('Bacongo', 'Brazzaville'),
    @('Djiri', 'Brazzaville'),
    @('Ile Mbamou', 'Brazzaville'),
    @('Madibou', 'Brazzaville'),
    @('Makélékélé', 'Brazzaville'),
    @('Mfilou', 'Brazzaville'),
    @('Moungali', 'Brazzaville'),
    @('Ouenzé', 'Brazzaville'),
    @('Poto-Poto', 'Brazzaville'),
    @('Talangai', 'Brazzaville'),
    @('Mossaka', 'Cuvette'),
    @('Owando', 'Cuvette'),
    @('Oyo-Alima', 'Cuvette'),
    @('Hinda', 'Kouilou'),
    @('Madingo-Kayes-Nzambi', 'Kouilou'),
    @('Mvouti-Kakamoeka', 'Kouilou'),
    @('Dolisie', 'Niari'),
    @('Kibangou', 'Niari'),
    @('Kimongo', 'Niari'),
    @('Mayoko', 'Niari'),
    @('Mossendjo', 'Niari'),
    @('Abala', 'Plateaux'),
    @('Djambala', 'Plateaux'),
    @('Gamboma', 'Plateaux'),
    @('Ngo-Mpouya', 'Plateaux'),
    @('Boko', 'Pool'),
    @('Goma Tsé-Tsé', 'Pool'),
    @('Igné-Ngabe-Mayama', 'Pool'),
    @('Kindamba', 'Pool'),
    @('Kinkala', 'Pool'),
    @('Kintele', 'Pool'),
    @('Mindouli', 'Pool')
)

# Clear the old region_list / district_list rows (15-34) before rewriting
# them, since the new lists are longer than the old ones.
$choices.Range("A15:D34").ClearContents()

# Write the region_list rows (15-21).
for ($i = 0; $i -lt $regions.Length; $i++) {
    $r = 15 + $i
    $choices.Cells.Item($r, 1).Value = "region_list"
    $choices.Cells.Item($r, 2).Value = $regions[$i]
    $choices.Cells.Item($r, 3).Value = $regions[$i]
}

# Write the district_list rows (22-58).
for ($i = 0; $i -lt $districts.Length; $i++) {
    $r = 22 + $i
    $pair = $districts[$i]
    $choices.Cells.Item($r, 1).Value = "district_list"
    $choices.Cells.Item($r, 2).Value = $pair[0]
    $choices.Cells.Item($r, 3).Value = $pair[0]
    $choices.Cells.Item($r, 4).Value = $pair[1]
}

# --- Make "choices" the active/selected sheet (was "survey") -------------
$choices.Activate()
$choices.Range("A38:A58").Select()
